$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers need to be
# forced to Text format first, otherwise Excel auto-converts them to numeric
# values (changing the cell type from string to number).
$ws.Range("D2").Value = "26.667.03"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "1.598.10"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.50"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.514"
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +0.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.56"
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "1.822.24"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "1.598.64"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.16"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").Value = "26.653.84"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "209.65"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("E21").Value = "  +3.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.30"
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("E23").Value = "  +1.71%  "
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.29"
$ws.Range("E25").Value = "  -1.44%  "
$ws.Range("E27").Value = "  -0.67%  "
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0515"
$ws.Range("E30").Value = "  +2.41%  "
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("D34").Value = "1.288.75"
$ws.Range("E34").Value = "  -0.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.619"
$ws.Range("E35").Value = "  -6.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.46"
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("E40").Value = "  +19.72%  "
$ws.Range("E41").Value = "  +2.40%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.785"
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.52"
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("D45").Value = "1.735.55"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.80"
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("E47").Value = "  -3.26%  "
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("E49").Value = "  +1.45%  "
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("E51").Value = "  +0.10%  "
